$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44249
$ws.Cells.Item(2, 12).Value = "Especial"
$ws.Cells.Item(2, 13).Value = 200
$ws.Cells.Item(2, 14).Value = 6000
$ws.Cells.Item(2, 15).Value = 7000
$ws.Cells.Item(2, 16).Value = 6500
$ws.Cells.Item(2, 19).Value = 2167

$ws.Cells.Item(3, 4).Value = 44249
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 160
$ws.Cells.Item(3, 14).Value = 4500
$ws.Cells.Item(3, 15).Value = 5000
$ws.Cells.Item(3, 16).Value = 4750
$ws.Cells.Item(3, 19).Value = 1583

$ws.Cells.Item(4, 4).Value = 44200
$ws.Cells.Item(4, 12).Value = "Especial"
$ws.Cells.Item(4, 13).Value = 50
$ws.Cells.Item(4, 14).Value = 4500
$ws.Cells.Item(4, 15).Value = 5000
$ws.Cells.Item(4, 16).Value = 4750
$ws.Cells.Item(4, 19).Value = 1583

$ws.Cells.Item(5, 4).Value = 44200
$ws.Cells.Item(5, 12).Value = "Primera"
$ws.Cells.Item(5, 13).Value = 80
$ws.Cells.Item(5, 14).Value = 3500
$ws.Cells.Item(5, 15).Value = 4000
$ws.Cells.Item(5, 16).Value = 3750
$ws.Cells.Item(5, 19).Value = 1250

$ws.Cells.Item(6, 4).Value = 44200
$ws.Cells.Item(6, 12).Value = "Segunda"
$ws.Cells.Item(6, 13).Value = 120
$ws.Cells.Item(6, 14).Value = 2500
$ws.Cells.Item(6, 15).Value = 3000
$ws.Cells.Item(6, 16).Value = 2750
$ws.Cells.Item(6, 19).Value = 917

$ws.Cells.Item(7, 4).Value = 44322
$ws.Cells.Item(7, 12).Value = "Especial"
$ws.Cells.Item(7, 13).Value = 200
$ws.Cells.Item(7, 14).Value = 7000
$ws.Cells.Item(7, 15).Value = 7500
$ws.Cells.Item(7, 16).Value = 7250
$ws.Cells.Item(7, 19).Value = 2417

$ws.Cells.Item(8, 4).Value = 44322
$ws.Cells.Item(8, 12).Value = "Primera"
$ws.Cells.Item(8, 13).Value = 160
$ws.Cells.Item(8, 14).Value = 6000
$ws.Cells.Item(8, 15).Value = 6500
$ws.Cells.Item(8, 16).Value = 6250
$ws.Cells.Item(8, 19).Value = 2083

$ws.Cells.Item(9, 4).Value = 44322
$ws.Cells.Item(9, 12).Value = "Segunda"
$ws.Cells.Item(9, 13).Value = 100
$ws.Cells.Item(9, 14).Value = 5000
$ws.Cells.Item(9, 15).Value = 5500
$ws.Cells.Item(9, 16).Value = 5250
$ws.Cells.Item(9, 19).Value = 1750

$ws.Cells.Item(10, 4).Value = 44351
$ws.Cells.Item(10, 12).Value = "Especial"
$ws.Cells.Item(10, 13).Value = 160
$ws.Cells.Item(10, 14).Value = 7500
$ws.Cells.Item(10, 15).Value = 8000
$ws.Cells.Item(10, 16).Value = 7750
$ws.Cells.Item(10, 19).Value = 2583

$ws.Cells.Item(11, 4).Value = 44351
$ws.Cells.Item(11, 12).Value = "Primera"
$ws.Cells.Item(11, 13).Value = 100
$ws.Cells.Item(11, 14).Value = 6000
$ws.Cells.Item(11, 15).Value = 6500
$ws.Cells.Item(11, 16).Value = 6250
$ws.Cells.Item(11, 19).Value = 2083

$ws.Cells.Item(12, 4).Value = 44351
$ws.Cells.Item(12, 12).Value = "Segunda"
$ws.Cells.Item(12, 13).Value = 200
$ws.Cells.Item(12, 14).Value = 4500
$ws.Cells.Item(12, 15).Value = 5000
$ws.Cells.Item(12, 16).Value = 4750
$ws.Cells.Item(12, 19).Value = 1583

$ws.Cells.Item(13, 4).Value = 44242
$ws.Cells.Item(13, 12).Value = "Especial"
$ws.Cells.Item(13, 13).Value = 50
$ws.Cells.Item(13, 14).Value = 7000
$ws.Cells.Item(13, 15).Value = 8000
$ws.Cells.Item(13, 16).Value = 7500
$ws.Cells.Item(13, 19).Value = 2500

$ws.Cells.Item(14, 4).Value = 44242
$ws.Cells.Item(14, 12).Value = "Primera"
$ws.Cells.Item(14, 13).Value = 90
$ws.Cells.Item(14, 14).Value = 6000
$ws.Cells.Item(14, 15).Value = 7000
$ws.Cells.Item(14, 16).Value = 6500
$ws.Cells.Item(14, 19).Value = 2167

$ws.Cells.Item(15, 4).Value = 44242
$ws.Cells.Item(15, 12).Value = "Segunda"
$ws.Cells.Item(15, 13).Value = 100
$ws.Cells.Item(15, 14).Value = 4000
$ws.Cells.Item(15, 15).Value = 5000
$ws.Cells.Item(15, 16).Value = 4500
$ws.Cells.Item(15, 19).Value = 1500

$ws.Cells.Item(16, 4).Value = 44334
$ws.Cells.Item(16, 12).Value = "Especial"
$ws.Cells.Item(16, 13).Value = 100
$ws.Cells.Item(16, 14).Value = 7000
$ws.Cells.Item(16, 15).Value = 8000
$ws.Cells.Item(16, 16).Value = 7500
$ws.Cells.Item(16, 19).Value = 2500

$ws.Cells.Item(17, 4).Value = 44334
$ws.Cells.Item(17, 12).Value = "Primera"
$ws.Cells.Item(17, 13).Value = 160
$ws.Cells.Item(17, 14).Value = 6000
$ws.Cells.Item(17, 15).Value = 7000
$ws.Cells.Item(17, 16).Value = 6500
$ws.Cells.Item(17, 19).Value = 2167

$ws.Cells.Item(18, 4).Value = 44334
$ws.Cells.Item(18, 12).Value = "Segunda"
$ws.Cells.Item(18, 13).Value = 120
$ws.Cells.Item(18, 14).Value = 6000
$ws.Cells.Item(18, 15).Value = 7000
$ws.Cells.Item(18, 16).Value = 6500
$ws.Cells.Item(18, 19).Value = 2167

$ws.Cells.Item(19, 4).Value = 44334
$ws.Cells.Item(19, 12).Value = "Tercera"
$ws.Cells.Item(19, 13).Value = 70
$ws.Cells.Item(19, 14).Value = 3500
$ws.Cells.Item(19, 15).Value = 4000
$ws.Cells.Item(19, 16).Value = 3750
$ws.Cells.Item(19, 19).Value = 1250

$ws.Cells.Item(20, 4).Value = 44172
$ws.Cells.Item(20, 12).Value = "Especial"
$ws.Cells.Item(20, 13).Value = 100
$ws.Cells.Item(20, 14).Value = 6500
$ws.Cells.Item(20, 15).Value = 7000
$ws.Cells.Item(20, 16).Value = 6750
$ws.Cells.Item(20, 19).Value = 2250

$ws.Cells.Item(21, 4).Value = 44172
$ws.Cells.Item(21, 12).Value = "Primera"
$ws.Cells.Item(21, 13).Value = 160
$ws.Cells.Item(21, 14).Value = 5500
$ws.Cells.Item(21, 15).Value = 6000
$ws.Cells.Item(21, 16).Value = 5750
$ws.Cells.Item(21, 19).Value = 1917

$ws.Cells.Item(22, 4).Value = 44172
$ws.Cells.Item(22, 12).Value = "Segunda"
$ws.Cells.Item(22, 13).Value = 160
$ws.Cells.Item(22, 14).Value = 5000
$ws.Cells.Item(22, 15).Value = 5500
$ws.Cells.Item(22, 16).Value = 5250
$ws.Cells.Item(22, 19).Value = 1750

$ws.Cells.Item(23, 4).Value = 44172
$ws.Cells.Item(23, 12).Value = "Tercera"
$ws.Cells.Item(23, 13).Value = 140
$ws.Cells.Item(23, 14).Value = 3500
$ws.Cells.Item(23, 15).Value = 4000
$ws.Cells.Item(23, 16).Value = 3750
$ws.Cells.Item(23, 19).Value = 1250

$ws.Cells.Item(24, 4).Value = 44389
$ws.Cells.Item(24, 12).Value = "Especial"
$ws.Cells.Item(24, 13).Value = 100
$ws.Cells.Item(24, 14).Value = 7500
$ws.Cells.Item(24, 15).Value = 8000
$ws.Cells.Item(24, 16).Value = 7750
$ws.Cells.Item(24, 19).Value = 2583

$ws.Cells.Item(25, 4).Value = 44389
$ws.Cells.Item(25, 12).Value = "Primera"
$ws.Cells.Item(25, 13).Value = 160
$ws.Cells.Item(25, 14).Value = 6000
$ws.Cells.Item(25, 15).Value = 7000
$ws.Cells.Item(25, 16).Value = 6500
$ws.Cells.Item(25, 19).Value = 2167

$ws.Cells.Item(26, 4).Value = 44389
$ws.Cells.Item(26, 12).Value = "Segunda"
$ws.Cells.Item(26, 13).Value = 200
$ws.Cells.Item(26, 14).Value = 5500
$ws.Cells.Item(26, 15).Value = 6000
$ws.Cells.Item(26, 16).Value = 5750
$ws.Cells.Item(26, 19).Value = 1917
